$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 31: Machine 30
$ws.Cells.Item(31, 1).Value = 10030
$ws.Cells.Item(31, 2).Value = "Machine 30"
$ws.Cells.Item(31, 3).Value = "70-5A-0F-8C-01-39"
$ws.Cells.Item(31, 4).Value = "FB5962911663"
$ws.Cells.Item(31, 5).Value = "192.168.0.356"
$ws.Cells.Item(31, 6).Value = 1001
$ws.Cells.Item(31, 7).Value = "eng"
$ws.Cells.Item(31, 8).Value = $true
$ws.Cells.Item(31, 9).Value = "superadmin"
$ws.Cells.Item(31, 10).Value = "now()"
$ws.Cells.Item(31, 11).Value = "now()"

# New row 32: Machine 31
$ws.Cells.Item(32, 1).Value = 10031
$ws.Cells.Item(32, 2).Value = "Machine 31"
$ws.Cells.Item(32, 4).Value = "FB5962911663"
$ws.Cells.Item(32, 5).Value = "192.168.0.357"
$ws.Cells.Item(32, 3).Value = "58-20-B1-DA-F3-FB"
$ws.Cells.Item(32, 6).Value = 1001
$ws.Cells.Item(32, 7).Value = "eng"
$ws.Cells.Item(32, 8).Value = $true
$ws.Cells.Item(32, 9).Value = "superadmin"
$ws.Cells.Item(32, 10).Value = "now()"
$ws.Cells.Item(32, 11).Value = "now()"

$ws.Range("L1:XFD1048576").Select() | Out-Null
